$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("6:6").Insert()
$ws.Range("A6:XFD6").Select()
